# TC06_Canine_Filter_FileFormat-tif.xlsx
# Fixed Diagnosis, FileAssociation, FileFormat, FileType, NeuteredStatus, PrimeDiseaseSite
#
# The "CasesTab" Neo4j query in cell B2 (startup sheet) had its trailing
# `Cohort` column removed from the RETURN clause.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
  MATCH (f:file)-[*]->(c)
    WHERE f.file_format IN ["tif"]  
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

# Strip the single trailing newline the here-string literal adds after the
# last line so the stored string matches the original text exactly.
$newQuery = $newQuery.TrimEnd("`r", "`n")

$ws.Range("B2").Value = $newQuery

# Move the selection back onto the edited cell (matches the saved view state).
$ws.Range("B2").Select()
